$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Constanze Shuruq Shuruq Bushnaq Josting -> split name update + birthdate shift
$ws.Range("C10").Value = "Constanze Shuruq Magnolia"
$ws.Range("D10").Value = "Josting"
$ws.Range("E10").Value = 36756

# Row 5: Abier Bushnaq birthdate shift (+1 year)
$ws.Range("E5").Value = 24222

# Row 8: Lucius Nabil Bushnaq birthdate shift (+2 years)
$ws.Range("E8").Value = 35129

# Update the active selection to B9 to match the saved view state
$ws.Range("B9").Select()
